# The edit swaps the data of row 12 and row 13 (all the species-observation
# fields), while the shared context columns (P, Q/R coordinates aside,
# S..W, Y, AA, AD, AE, AG, AT, AW, AY) stay associated with the same row
# index. In effect, the "Tjäder" record (with its K/L/M/N activity data)
# moves from row 12 to row 13, and the "Garnlav" record moves from row 13
# to row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 becomes the former row-13 ("Garnlav") record ---
$ws.Range("A12").Value = 130981935
$ws.Range("B12").Value = 79243
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("Q12").Value = 437656
$ws.Range("R12").Value = 6792404
$ws.Range("AX12").Value = "Eva Löfqvist"

# --- Row 13 becomes the former row-12 ("Tjäder") record ---
$ws.Range("A13").Value = 130981909
$ws.Range("B13").Value = 57073
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 100138
$ws.Range("F13").Value = "Tjäder"
$ws.Range("G13").Value = "Tetrao urogallus"
$ws.Range("H13").Value = "Linnaeus, 1758"
# K13/L13/N13 need to exist as blank (but present) cells, matching the
# blank I13 cell already on the sheet, so copy that blank cell into them
# instead of assigning "" (which would remove the cell altogether).
$ws.Range("I13").Copy($ws.Range("K13"))
$ws.Range("I13").Copy($ws.Range("L13"))
$ws.Range("M13").Value = "färsk spillning"
$ws.Range("I13").Copy($ws.Range("N13"))
$ws.Range("Q13").Value = 437657
$ws.Range("R13").Value = 6792398
$ws.Range("AX13").Value = "Eva Löfqvist, Alfhild Sehlin"
